$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  ,@("M-SL 営業",45209,$null,45271,"Expired",157,160,0.98125)
  ,@("M-MN 管理",45209,$null,45271,"Expired",68.75,80,0.859375)
  ,@("M-BT 出張",45209,$null,45271,"Expired",40,80,0.5)
  ,@("M-MT ミーティング",45209,$null,45271,"Expired",68,160,0.425)
  ,@("PP-IC-IN インストール",45224,45226,45271,"Expired",2,4,0.5)
  ,@("PP-IC-OP オペレーション教育",45224,45226,45271,"Expired",-14.5,24,-0.6041666666666666)
  ,@("PT-OC OrCAD",45224,45224,45271,"Expired",3,8,0.375)
  ,@("PP-MC メカ全体",45224,45230,45271,"Expired",-1.5,16,-0.09375)
  ,@("PP-HD-BA 基本機構の理解",45224,45230,45271,"Expired",-2,16,-0.125)
  ,@("PP-HD-MD 対象メディアの理解",45224,45225,45271,"Expired",7,8,0.875)
  ,@("PP-EL-MP メイン基板",45224,45225,45271,"Expired",-10.5,16,-0.65625)
  ,@("PP-EL-CP 子基板",45224,45225,45271,"Expired",11,16,0.6875)
  ,@("PP-EL-PP 電源",45224,45225,45271,"Expired",5,8,0.625)
  ,@("PP-SW-FW ファーム",45224,45225,45271,"Expired",30,40,0.75)
  ,@("PP-IS-HD ヘッド",45230,$null,45271,"Expired",9,16,0.5625)
  ,@("PI-EV-SV サービス性",45230,$null,45271,"Expired",9.75,16,0.609375)
  ,@("PI-EV-IV 固有評価",45230,$null,45271,"Expired",0,16,0)
  ,@("PI-EV-MS MES評価",45230,$null,45271,"Expired",-1,8,-0.125)
  ,@("PI-EV-BE 基板等の電気的評価",45230,$null,45271,"Expired",0,8,0)
  ,@("PI-EV-PE 製品の電気的評価",45230,$null,45271,"Expired",0,8,0)
  ,@("PI-EV-SS 安全規格の評価",45230,$null,45271,"Expired",0,8,0)
  ,@("PI-EV-GN 全体動作の評価",45230,$null,45271,"Expired",0,8,0)
)

$startRow = 2
$r = $startRow
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]

  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"

  if ($row[2] -eq $null) {
    $ws.Cells.Item($r, 3).ClearContents()
  } else {
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
  }

  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD"

  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]

  $r = $r + 1
}

Write-Output "Updated rows 2 through $($r - 1)"
